# edit.ps1 - applies the "output generated at 456a3b4" update to the workbook
$wb = $excel.ActiveWorkbook

# ===================== Sheet: 展览 (Exhibition) =====================
$ws1 = $wb.Worksheets.Item("展览")

# "想去人数" (want-to-go count) bumps
$ws1.Cells.Item(4, 6).Value = 5086
$ws1.Cells.Item(5, 6).Value = 5086
$ws1.Cells.Item(6, 6).Value = 119
$ws1.Cells.Item(7, 6).Value = 152
$ws1.Cells.Item(11, 6).Value = 170
$ws1.Cells.Item(12, 6).Value = 8406
$ws1.Cells.Item(13, 6).Value = 8406
$ws1.Cells.Item(17, 6).Value = 610
$ws1.Cells.Item(18, 6).Value = 2532
$ws1.Cells.Item(24, 6).Value = 2524
$ws1.Cells.Item(27, 6).Value = 6408
$ws1.Cells.Item(28, 6).Value = 187
$ws1.Cells.Item(33, 6).Value = 6901
$ws1.Cells.Item(35, 6).Value = 32
$ws1.Cells.Item(38, 6).Value = 10
$ws1.Cells.Item(40, 6).Value = 27
$ws1.Cells.Item(43, 6).Value = 2531
$ws1.Cells.Item(45, 6).Value = 69
$ws1.Cells.Item(46, 6).Value = 1128
$ws1.Cells.Item(48, 6).Value = 517
$ws1.Cells.Item(49, 6).Value = 2220
$ws1.Cells.Item(50, 6).Value = 76

# Row 19 content shifts up to what was row 20 (模型博览会), with refreshed counts
$ws1.Cells.Item(19, 2).Value = "'2024-04-19"
$ws1.Cells.Item(19, 3).Value = "北京·第22届中国国际模型博览会"
$ws1.Cells.Item(19, 4).Value = "北京展览馆 北京展览馆"
$ws1.Cells.Item(19, 5).Value = "2024.04.19 10:00-04.21 17:00"
$ws1.Cells.Item(19, 6).Value = 6327
$ws1.Cells.Item(19, 7).Value = 13.5
$ws1.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82425"
$ws1.Cells.Item(19, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/9nkCFSHm1709710888611.jpeg"

# Row 20 content shifts up to what was row 21 (QMQ动漫游戏嘉年华), with refreshed counts
$ws1.Cells.Item(20, 2).Value = "'2024-04-20"
$ws1.Cells.Item(20, 3).Value = "北京·QMQ动漫游戏嘉年华"
$ws1.Cells.Item(20, 4).Value = "小关路39号 北投购物公园"
$ws1.Cells.Item(20, 5).Value = "2024.04.20 10:00-04.21 17:00"
$ws1.Cells.Item(20, 6).Value = 2303
$ws1.Cells.Item(20, 7).Value = 63
$ws1.Cells.Item(20, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81982"
$ws1.Cells.Item(20, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/lyPb1fLO1708569465126.jpeg"

# Row 21 content shifts up to what was row 22 (亚力传感器走秀派对); date (col B) unchanged
$ws1.Cells.Item(21, 3).Value = "北京·亚力传感器走秀派对"
$ws1.Cells.Item(21, 4).Value = "旧鼓楼大街51号(鼓楼大街地铁站G东南口步行250米) MOONEE 暮霓"
$ws1.Cells.Item(21, 5).Value = "2024.04.20 19:00-04.21 02:00"
$ws1.Cells.Item(21, 6).Value = 6
$ws1.Cells.Item(21, 7).Value = 68
$ws1.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83540"
$ws1.Cells.Item(21, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/ZDsD1X9t1711523212670.jpeg"

# Row 22 becomes a brand-new entry (国际电影节·光影未来...); date (col B) unchanged
$ws1.Cells.Item(22, 3).Value = "北京·国际电影节·光影未来“万游影力”影游动画狂欢派对"
$ws1.Cells.Item(22, 4).Value = "半截塔路53号首创郎园station西门 郎园station中央车站文化中心"
$ws1.Cells.Item(22, 5).Value = "2024.04.20 10:00-04.21 17:00"
$ws1.Cells.Item(22, 6).Value = 1
$ws1.Cells.Item(22, 7).Value = 75
$ws1.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83564"
$ws1.Cells.Item(22, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/yHtZ4o5y1711707895213.jpeg"

# ===================== Sheet: 演出 (Performance) =====================
$ws2 = $wb.Worksheets.Item("演出")

# "想去人数" (want-to-go count) bumps
$ws2.Cells.Item(2, 6).Value = 14
$ws2.Cells.Item(3, 6).Value = 167
$ws2.Cells.Item(5, 6).Value = 57
$ws2.Cells.Item(6, 6).Value = 15
$ws2.Cells.Item(12, 6).Value = 155
$ws2.Cells.Item(13, 6).Value = 9

# ===================== Sheet: 全部类型 (All Types) =====================
$ws4 = $wb.Worksheets.Item("全部类型")

# "想去人数" (want-to-go count) bumps
$ws4.Cells.Item(3, 6).Value = 5086
$ws4.Cells.Item(4, 6).Value = 5086
$ws4.Cells.Item(5, 6).Value = 119
$ws4.Cells.Item(6, 6).Value = 152
$ws4.Cells.Item(10, 6).Value = 170
$ws4.Cells.Item(11, 6).Value = 8406
$ws4.Cells.Item(12, 6).Value = 8406
$ws4.Cells.Item(15, 6).Value = 610
$ws4.Cells.Item(16, 6).Value = 2532
$ws4.Cells.Item(17, 6).Value = 167
$ws4.Cells.Item(18, 6).Value = 6327
$ws4.Cells.Item(19, 6).Value = 2303
$ws4.Cells.Item(20, 6).Value = 57
$ws4.Cells.Item(22, 6).Value = 2524
$ws4.Cells.Item(24, 6).Value = 15
$ws4.Cells.Item(27, 6).Value = 6408
$ws4.Cells.Item(28, 6).Value = 187
$ws4.Cells.Item(33, 6).Value = 6901
$ws4.Cells.Item(35, 6).Value = 32
$ws4.Cells.Item(37, 6).Value = 10
$ws4.Cells.Item(39, 6).Value = 27
$ws4.Cells.Item(42, 6).Value = 2531
$ws4.Cells.Item(44, 6).Value = 69
$ws4.Cells.Item(45, 6).Value = 1128
$ws4.Cells.Item(47, 6).Value = 517
$ws4.Cells.Item(48, 6).Value = 155
$ws4.Cells.Item(49, 6).Value = 2220
$ws4.Cells.Item(50, 6).Value = 76
